$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 224, shifting existing rows 224:324 down by one.
$ws.Rows.Item(224).Insert()

# Fill in the new row's data (matches the row that used to be at 224, but with
# updated price/date observations for the newly added weekly entry).
$ws.Cells.Item(224, 1).Value = 9
$ws.Cells.Item(224, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(224, 3).Value = "Metropolitana"
$ws.Cells.Item(224, 4).Value = 44825
$ws.Cells.Item(224, 5).Value = 13
$ws.Cells.Item(224, 6).Value = 300000001
$ws.Cells.Item(224, 7).Value = "Rabanito"
$ws.Cells.Item(224, 8).Value = "Sin especificar"
$ws.Cells.Item(224, 9).Value = "Primera"
$ws.Cells.Item(224, 10).Value = 7900
$ws.Cells.Item(224, 11).Value = 2500
$ws.Cells.Item(224, 12).Value = 3000
$ws.Cells.Item(224, 13).Value = 2750
$ws.Cells.Item(224, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(224, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(224, 16).Value = 28
$ws.Cells.Item(224, 17).Value = 100
$ws.Cells.Item(224, 18).Value = "Hortaliza"
